$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17821729183197
$ws.Range("B1").Value = 2.359478712081909
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.070424556732178
$ws.Range("E1").Value = 1.034740924835205
